# doc: end Alpha 1 and update Alpha 2 phases
$wb = $excel.ActiveWorkbook

$wsAlpha1 = $wb.Worksheets.Item("ALPHA 1")
$wsAlpha2 = $wb.Worksheets.Item("ALPHA 2")

# --- ALPHA 1: close out the phase ---------------------------------------
# Task "Gestione database per Customers" note is resolved -> clear the note
$wsAlpha1.Range("E9").ClearContents() | Out-Null

# Last task ("Autorecognition serial") is now Done -> mark the Done column
$wsAlpha1.Range("F15").Copy() | Out-Null
$wsAlpha1.Range("F16").PasteSpecial(-4122) | Out-Null

# Move the "currently open" selection off Alpha 1
$wsAlpha1.Range("E25").Select() | Out-Null

# --- ALPHA 2: refresh the phase's task list -------------------------------
$wsAlpha2.Range("C7").Value = "Simulazione di più box sullo stesso bridge"
$wsAlpha2.Range("E7").ClearContents() | Out-Null

$wsAlpha2.Range("C9").Value = "Separazione Applicativo/Basso livello"
$wsAlpha2.Range("D9").Value = "Simo"

$wsAlpha2.Range("C11").Value = "Sviluppo Camera e comunicazione con Bridge"
$wsAlpha2.Range("D11").Value = "Fre/Simo"
$wsAlpha2.Range("D11").Borders.LineStyle = -4142

$wsAlpha2.Range("C12").Value = "ClientServer(Bot Telegram?) per visualizzare log accessi?"
$wsAlpha2.Range("D12").Value = "Fre/Simo"

$wsAlpha2.Range("C15").Value = "Dubbi"
$wsAlpha2.Range("C15").WrapText = $false

$wsAlpha2.Range("C17").Value = "Central: due sketch: train + main?"
$wsAlpha2.Range("C18").Value = "Local: uno sketch?"

# ALPHA 2 becomes the active tab / sheet
$wsAlpha2.Activate() | Out-Null
$wsAlpha2.Range("C17").Select() | Out-Null
